$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cells to reflect the new terminology.
$ws.Range("D1").Value = "CapacityForecast"
$ws.Range("E1").Value = "EffortForecast"
$ws.Range("F1").Value = "CapacityDone"
$ws.Range("G1").Value = "EffortDone"

# Update the active selection as recorded in the sheet view.
$ws.Range("F2").Select()
